$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Encoder Left")
$ws2 = $wb.Worksheets.Item("Encoder Right")

# ---------------------------------------------------------------------------
# 1. Populate "Encoder Left" (sheet1) test-description row (row 2) and the
#    trial-results header row (row 3), then the trial data rows (4-7).
#    Cells are written in the same left-to-right-ish order the shared
#    string table in the target file implies (C2,E2,F2,H2,G2,B2,I2,J2 then
#    the row-3 headers) so new shared strings land at matching indices.
# ---------------------------------------------------------------------------

$ws1.Range("C2").Value = "32 turns per 15000 ticks"
$ws1.Range("E2").Value = "end of line in robot control"
$ws1.Range("F2").Value = "robot control model in routines subsystem"
$ws1.Range("H2").Value = "all devices connected except blue and green LEDs and bluetooth"
$ws1.Range("G2").Value = "motor connected normally in its position and a flag is attached to detect each round"
$ws1.Range("B2").Value = "give PWM for motor for 15000 ticks, and monitor how many turns it did. Repeat this sequence for different PWMs. 100, 150, 200 and 250"
$ws1.Range("I2").Value = "matlab model"
$ws1.Range("J2").Value = "leftMotorTicksPerCm = `nor`nleftMotorCmPerTick = "

$ws1.Range("B3").Value = "PWM"
$ws1.Range("C3").Value = "number of rounds"
$ws1.Range("D3").Value = "ticks/round"
$ws1.Range("E3").Value = "wheelSize (cm)"
$ws1.Range("F3").Value = "Ticks/Cm"

# Trial rows: PWM, number of rounds, ticks/round (formula), wheel size, Ticks/Cm (formula)
$ws1.Range("B4").Value = 100
$ws1.Range("B5").Value = 150
$ws1.Range("B6").Value = 200
$ws1.Range("B7").Value = 250

$ws1.Range("C4").Value = 32
$ws1.Range("C5").Value = 32
$ws1.Range("C6").Value = 32
$ws1.Range("C7").Value = 32

$ws1.Range("D4").Formula = "=15000/C4"
$ws1.Range("D5:D7").Formula = "=15000/C5"

$ws1.Range("E4").Value = 12
$ws1.Range("E5").Value = 12
$ws1.Range("E6").Value = 12
$ws1.Range("E7").Value = 12

$ws1.Range("F4").Formula = "=D4/(2*3.14*(E4/2))"
$ws1.Range("F5:F7").Formula = "=D5/(2*3.14*(E5/2))"

# ---------------------------------------------------------------------------
# 2. Row heights / column widths on "Encoder Left" for the now-wrapped rows.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).RowHeight = 100.8
$ws1.Rows.Item(3).RowHeight = 28.8

$ws1.Columns.Item(1).ColumnWidth = 5.666666666666667
$ws1.Columns.Item(10).ColumnWidth = 24

# ---------------------------------------------------------------------------
# 3. Duplicate the fully-populated "Encoder Left" layout into the still-empty
#    "Encoder Right" sheet (copy carries over values, formulas and styles).
# ---------------------------------------------------------------------------
$ws1.Range("A1:J9").Copy($ws2.Range("A1"))

# Encoder Right's description row points at the right-motor figures instead.
$ws2.Range("J2").Value = "rightMotorTicksPerCm = `nor`nrightMotorCmPerTick = "

# Row 8 on Encoder Right should stay blank in columns A-D, same as Encoder Left.
$ws2.Range("A8:D8").Clear()

# Row heights differ slightly between the two sheets.
$ws2.Rows.Item(2).RowHeight = 86.4
$ws2.Rows.Item(3).RowHeight = 28.8
$ws2.Rows.Item(9).RowHeight = 15

# Column widths specific to Encoder Right.
$ws2.Columns.Item(2).ColumnWidth = 22.109375
$ws2.Columns.Item(3).ColumnWidth = 12.21875
$ws2.Columns.Item(4).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(5).ColumnWidth = 14.88671875
$ws2.Columns.Item(6).ColumnWidth = 16.5
$ws2.Columns.Item(7).ColumnWidth = 25
$ws2.Columns.Item(8).ColumnWidth = 17
$ws2.Columns.Item(9).ColumnWidth = 16.88671875
$ws2.Columns.Item(10).ColumnWidth = 21

# ---------------------------------------------------------------------------
# 4. Selections: Encoder Left keeps a whole-table selection and is no longer
#    the active tab; Encoder Right becomes the active tab with J3 selected.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1:J9").Select()

$ws2.Activate()
$ws2.Range("J3").Select()

# ---------------------------------------------------------------------------
# 5. Page setup tweak recorded on Encoder Left.
# ---------------------------------------------------------------------------
$ws1.PageSetup.Orientation = 1

Write-Output "done"
